# Adds the new "stim details" block at the bottom of the sheet (rows 27-36),
# mirroring the header/table structure already used at the top of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 27 - section label
$ws.Range("A27").Value = "stim details"

# Row 28 - new table header row
$ws.Range("A28").Value = "month"
$ws.Range("B28").Value = "word_type"
$ws.Range("C28").Value = "need_audio"
$ws.Range("D28").Value = "need_image"
$ws.Range("E28").Value = "word"
$ws.Range("F28").Value = "count"
$ws.Range("G28").Value = "find images"

# Rows 29-32: video rows, months 6, 6, 7, 7
# Rows 33-36: audio rows, months 6, 6, 7, 7
$months = @(6, 6, 7, 7, 6, 6, 7, 7)
$types  = @("video", "video", "video", "video", "audio", "audio", "audio", "audio")

for ($i = 0; $i -lt $months.Length; $i++) {
    $row = 29 + $i
    $ws.Cells.Item($row, 1).Value = $months[$i]
    $ws.Cells.Item($row, 2).Value = $types[$i]
}
